# Update the "想去人数" (want-to-go count) figures on the "展览" and
# "全部类型" sheets to reflect the newly generated data snapshot.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 240
$ws1.Range("F6").Value = 13
$ws1.Range("F7").Value = 5792
$ws1.Range("F8").Value = 5129
$ws1.Range("F9").Value = 27
$ws1.Range("F10").Value = 56
$ws1.Range("F13").Value = 216
$ws1.Range("F14").Value = 26

# --- Sheet "全部类型" (all types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 240
$ws4.Range("F6").Value = 13
$ws4.Range("F7").Value = 5792
$ws4.Range("F8").Value = 5129
$ws4.Range("F9").Value = 27
$ws4.Range("F10").Value = 56
$ws4.Range("F13").Value = 216
$ws4.Range("F14").Value = 79
$ws4.Range("F16").Value = 26
